$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.405.81"
$ws.Range("E2").Value = "  +1.79%  "
$ws.Range("D3").Value = "3.940.74"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'490.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("D6").Value = "'146.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("D7").Value = "'0.624"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").Value = "'0.997"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "'0.736"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.10%  "
$ws.Range("D10").Value = "'0.177"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.53%  "
$ws.Range("D11").Value = "'0.0000347"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.07%  "
$ws.Range("D12").Value = "'42.96"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("D13").Value = "'10.48"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.47%  "
$ws.Range("D14").Value = "4.561.91"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").Value = "3.931.33"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("D16").Value = "'14.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.55%  "
$ws.Range("D18").Value = "'19.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.92%  "
$ws.Range("D19").Value = "'1.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.37%  "
$ws.Range("D20").Value = "69.224.81"
$ws.Range("E20").Value = "  +1.30%  "
$ws.Range("D21").Value = "'438.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.79%  "
$ws.Range("D22").Value = "'3.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.01%  "
$ws.Range("D23").Value = "'14.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.41%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'89.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.68%  "
$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").Value = "'12.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +12.88%  "
$ws.Range("D26").Value = "'3.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.46%  "
$ws.Range("D27").Value = "'11.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.31%  "
$ws.Range("D28").Value = "'37.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.27%  "
$ws.Range("D29").Value = "'5.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.58%  "
$ws.Range("D30").Value = "'709.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.90%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").Value = "'13.52"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.132"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.81%  "
$ws.Range("D33").Value = "'2.89"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.98%  "
$ws.Range("D34").Value = "'0.471"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +24.61%  "
$ws.Range("D35").Value = "0.0₃0912"
$ws.Range("E35").Value = "  -4.85%  "
$ws.Range("D36").Value = "'61.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.12%  "
$ws.Range("D37").Value = "'6.02"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.09%  "
$ws.Range("D38").Value = "'40.58"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.00%  "
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").Value = "  +2.14%  "
$ws.Range("D43").Value = "'2.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.24%  "
$ws.Range("E44").Value = "  -0.88%  "
$ws.Range("E45").Value = "  +2.54%  "
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").Value = "'3.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.60%  "
$ws.Range("D48").Value = "'3.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.56%  "
$ws.Range("D49").Value = "0.0₆0360"
$ws.Range("E49").Value = "  +11.22%  "
$ws.Range("D50").Value = "'3.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.19%  "
$ws.Range("D51").Value = "'2.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.40%  "
